$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C10 value from 18 to 100 (numeric)
$ws.Range("C10").Value = 100
